$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns L, M, N ---
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Match the header styling (bold, centered, bordered) used by the existing
# header cells by copying formats only from K1 (already styled) onto L1:N1.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)

# --- Data rows 2-7: new values for columns L, M, N ---
$ws.Range("L2").Value = 91.46496460658059
$ws.Range("M2").Value = 261310
$ws.Range("N2").Value = 319.4498777506112

$ws.Range("L3").Value = 72.47879000332269
$ws.Range("M3").Value = 2243
$ws.Range("N3").Value = 186.9166666666667

$ws.Range("L4").Value = 90.1751955003851
$ws.Range("M4").Value = 202038
$ws.Range("N4").Value = 147.5807158509861

$ws.Range("L5").Value = 89.1267344852567
$ws.Range("M5").Value = 1608
$ws.Range("N5").Value = 114.8571428571429

$ws.Range("L6").Value = 19.55931057734819
$ws.Range("M6").Value = 2208
$ws.Range("N6").Value = 14.52631578947368

$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
